# 31 October first commit
# The Katalon "AUTO_30OCT_TOOL_KAT_*" run tool keeps stamping cell A2 with a
# new timestamped token every time the suite executes; each run appends a
# fresh shared string and leaves A2 pointing at the newest one. Replay the
# same sequence of writes here so A2 ends up on the final 31-Oct token.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540891420271"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540891686239"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540891906051"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540892353385"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540892612327"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540892854267"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540893266411"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540893458319"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540893699787"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540893986430"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540894238982"
$ws.Range("A2").Value = "AUTO_30OCT_TOOL_KAT_1540894475278"
